$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$shriekmaw = "('Shriekmaw', ['{4}{B}', 'Creature — Elemental', 'Fear (This creature can" + [char]8217 + "t be blocked except by artifact creatures and/or black creatures.)', 'When Shriekmaw enters the battlefield, destroy target nonartifact, nonblack creature.', 'Evoke {1}{B} (You may cast this spell for its evoke cost. If you do, it" + [char]8217 + "s sacrificed when it enters the battlefield.)', '3/2'])"

$packmaster = "(`"Wren's Run Packmaster`", ['{3}{G}', 'Creature — Elf Warrior', 'Champion an Elf (When this creature enters the battlefield, sacrifice it unless you exile another Elf you control. When this creature leaves the battlefield, that card returns to the battlefield.)', '{2}{G}: Create a 2/2 green Wolf creature token.', 'Wolves you control have deathtouch.', '5/5'])"

$ws.Range("A2").Value = $shriekmaw
$ws.Range("A3").Value = $packmaster

# Delete rows 4 through 15 (old leftover data)
$ws.Range("A4:A15").EntireRow.Delete()
